# svm L1 L2 Separately
# Fill in the previously-empty SVM "L1/L2 separately" rows (59-65) in
# columns B:J with the newly measured bug-detection rates, and record the
# matching bug-description text in column M for rows 61-65 (rows 59 and 60
# already carried their M text from before). Column K (COUNTIF) and the
# row-68 totals are formula-driven and recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59
$ws.Range("B59").Value = 0.44
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 0.54
$ws.Range("E59").Value = 0.54
$ws.Range("F59").Value = 0.12
$ws.Range("G59").Value = 0.38
$ws.Range("H59").Value = 0.68
$ws.Range("I59").Value = 0.08
$ws.Range("J59").Value = 0.06

# Row 60
$ws.Range("B60").Value = 0
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 0.78
$ws.Range("E60").Value = 0.74
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0.08
$ws.Range("H60").Value = 0.04
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0.12

# Row 61
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 0.96
$ws.Range("E61").Value = 0.64
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0.02
$ws.Range("H61").Value = 0.04
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0.08
$ws.Range("M61").Value = "b_curr =b+0.12729193727342922 - step_size * gradient_b"

# Row 62
$ws.Range("B62").Value = 0
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 0.02
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0.16
$ws.Range("G62").Value = 0.56000000000000005
$ws.Range("H62").Value = 1
$ws.Range("I62").Value = 0.44
$ws.Range("J62").Value = 0.62
$ws.Range("M62").Value = "b_curr = b -step_size+0.5539443047647777 * gradient_b"

# Row 63
$ws.Range("B63").Value = 0
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 0.96
$ws.Range("E63").Value = 0.98
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0.38
$ws.Range("I63").Value = 0.26
$ws.Range("J63").Value = 0
$ws.Range("M63").Value = "b_curr = b - step_size *gradient_b+0.5293549766700935"

# Row 64
$ws.Range("B64").Value = 0
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0.02
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("M64").Value = " w = (1 -gamma-0.008003863215121487) * w_curr + gamma * w_prev"

# Row 65
$ws.Range("B65").Value = 0.02
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 0.94
$ws.Range("E65").Value = 0.92
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0.54
$ws.Range("H65").Value = 0.52
$ws.Range("I65").Value = 0.04
$ws.Range("J65").Value = 0.04
$ws.Range("M65").Value = "w = (1 - gamma) *w_curr-0.9117207568581369 + gamma * w_prev"

# Recalculate so K (COUNTIF) and the row 68 totals pick up the new data.
$excel.CalculateFull()

# Match the final selection/scroll position left by the edit session.
$ws.Range("A43").Select()
$ws.Range("J65").Select()
